$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data for rows 11-13 is cyclically rotated:
#   new row 11 <- old row 13
#   new row 12 <- old row 11
#   new row 13 <- old row 12
# Capture the "old" values for the affected columns before overwriting.

$cols = @("A","B","E","F","G","H","Q","R","Z","AB")

$old11 = @{}
$old12 = @{}
$old13 = @{}
foreach ($col in $cols) {
    $old11[$col] = $ws.Range($col + "11").Value2
    $old12[$col] = $ws.Range($col + "12").Value2
    $old13[$col] = $ws.Range($col + "13").Value2
}

foreach ($col in $cols) {
    $ws.Range($col + "11").Value = $old13[$col]
    $ws.Range($col + "12").Value = $old11[$col]
    $ws.Range($col + "13").Value = $old12[$col]
}

# H11 must end up empty (row 13's H was empty) and H13 must be populated
# (row 12's H value). Since old H13 was empty (Value2 -> $false/empty),
# explicitly clear H11 to guarantee an empty cell rather than a literal
# "False"/0 artifact from an empty Variant.
if (-not $old13["H"]) {
    $ws.Range("H11").ClearContents()
}
